$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Beta)
$ws.Range("F2").Value = 34.2546190917304
$ws.Range("G2").Value = 33.65494406955679
$ws.Range("H2").Value = 34.87597401633013
$ws.Range("I2").Value = 3.827526590341443
$ws.Range("J2").Value = 3.799985332062707
$ws.Range("K2").Value = 3.854913529396547
$ws.Range("L2").Value = 0.2683947526569028
$ws.Range("M2").Value = 0.2663134746221297
$ws.Range("N2").Value = 0.2705001393067226

# Row 3 (Gamma)
$ws.Range("F3").Value = 0.0002453897857124177
$ws.Range("G3").Value = 0.0000000952519424010228
$ws.Range("H3").Value = 0.0006785031254740325
$ws.Range("I3").Value = 0.0002305434406631316
$ws.Range("J3").Value = 0.00000009012529177433422
$ws.Range("K3").Value = 0.0006375236797138979
$ws.Range("L3").Value = 0.0002459749475615888
$ws.Range("M3").Value = 0.00000009605756640378798
$ws.Range("N3").Value = 0.0006803661664941629

# Row 4 (Beta + Gamma)
$ws.Range("F4").Value = 34.25486448151611
$ws.Range("G4").Value = 33.65494416480874
$ws.Range("H4").Value = 34.87665251945561
$ws.Range("I4").Value = 3.827757133782107
$ws.Range("J4").Value = 3.799985422187999
$ws.Range("K4").Value = 3.855551053076261
$ws.Range("L4").Value = 0.2686407276044644
$ws.Range("M4").Value = 0.2663135706796961
$ws.Range("N4").Value = 0.2711805054732168
